# Fixed import of CDP to include cost of investment
# Insert a new "Cost Of Investment *" column between Gross (C) and Reinvestment (D),
# shifting the existing Reinvestment/Date/Generate Payments/Payments Paid columns
# one column to the right, and populate the new column with the cost-of-investment
# figures for each distribution row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns D:G one column to the right to make room for the new column.
$ws.Range("D:D").Insert()

# New header for the inserted column.
$ws.Range("D1").Value = "Cost Of Investment *"

# New values for the inserted column (cost of investment for each row).
$ws.Range("D2").Value = 800000
$ws.Range("D3").Value = 1500000
$ws.Range("D4").Value = 2500000

# Match the numeric formatting used by the neighbouring Gross/Reinvestment columns.
$ws.Range("D2:D4").NumberFormat = $ws.Range("E2:E4").NumberFormat

# Widen the new column like the source workbook (~17.625, matching the Gross column's
# width). The host's ColumnWidth setter quantizes to whole pixels, so 16.9 is the input
# that lands closest to the target 17.625 "characters" width.
$ws.Range("D:D").ColumnWidth = 16.9

# Move the active selection to D5 to match the post-edit cursor position.
$null = $ws.Range("D5").Select()
